$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RBFF")

# Biomass (column G) no longer shifts fully to electricity (row 2);
# it now stays as biomass (row 7) -- i.e. no fuel shift for biomass.
$ws.Range("G2").Value = 0
$ws.Range("G7").Value = 1

# Make RBFF the active/visible sheet, matching the saved view state.
$ws.Activate()
